# Update two weeks of ARS testing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old ad-hoc "Länge für Einheit / Zahl" helper calc block that
#    used to sit next to week 5's data (rows 169-175, columns I:J). Clear()
#    removes both content and formatting so the cells disappear entirely.
#    (We grab the two highlighted-style rows first so we can stamp the same
#    formats onto the rebuilt block further down.)
# ---------------------------------------------------------------------------
$ws.Range("I174:J174").Copy()
$ws.Range("Z1:AA1").PasteSpecial(-4122)
$ws.Range("I175:J175").Copy()
$ws.Range("Z2:AA2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I169:J175").Clear()

# ---------------------------------------------------------------------------
# 2. Corrected measurement for week 7 / age group 5-14 (row 185, col E).
# ---------------------------------------------------------------------------
$ws.Range("E185").Value = 9.7

# ---------------------------------------------------------------------------
# 3. Re-create the helper calc block, now attached to week 7's data
#    (rows 184-190, columns I:J), ending with a new "Source: 2021-03-23"
#    separator (B190) that introduces week 8.
# ---------------------------------------------------------------------------
$ws.Range("I184").Value = "Länge für Einheit (cm)"
$ws.Range("J184").Value = 16.75

$ws.Range("I185").Value = "#:"
$ws.Range("J185").Value = 10

$ws.Range("I186").Value = "# / cm"
$ws.Range("J186").Formula = "=J185 / J184"

$ws.Range("I187").Value = "Achsenabschnitt (cm)"
$ws.Range("J187").Value = 0

$ws.Range("I188").Value = "Achsenabschnitt (#)"
$ws.Range("J188").Value = 5

$ws.Range("I189").Value = "Gemessene Höhe (cm)"
$ws.Range("J189").Value = 7.45

$ws.Range("B190").Value = "Source: 2021-03-23"
$ws.Range("I190").Value = "Zahl"
$ws.Range("J190").Formula = "=(J189-J187)*J186+J188"

# Stamp the highlighted ("measured value") and bold ("computed value") formats
# saved off above onto the new block, same relative rows:
#   old I174/J174 (bold red)  -> new I189/J189
#   old I175/J175 (bold)      -> new I190/J190
$ws.Range("Z1:AA1").Copy()
$ws.Range("I189:J189").PasteSpecial(-4122)
$ws.Range("Z2:AA2").Copy()
$ws.Range("I190:J190").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("Z1:AA2").Clear()

# ---------------------------------------------------------------------------
# 4. New raw testing data, week 8 (rows 191-196).
# ---------------------------------------------------------------------------
$wk8 = @(
    @(2021, 8, "0-4",   9333,   6.7),
    @(2021, 8, "5-14",  14000,  9.8),
    @(2021, 8, "15-34", 96666,  6.6),
    @(2021, 8, "35-59", 154000, 5.9),
    @(2021, 8, "60-79", 86667,  5.7),
    @(2021, 8, ">=80",  45000,  6.6)
)
$r = 191
foreach ($row in $wk8) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5. Separator introducing week 9 (row 197) - same fetch date as week 8.
# ---------------------------------------------------------------------------
$ws.Range("B197").Value = "Source: 2021-03-23"

# ---------------------------------------------------------------------------
# 6. New raw testing data, week 9 (rows 198-203).
# ---------------------------------------------------------------------------
$wk9 = @(
    @(2021, 9, "0-4",   14000,  6),
    @(2021, 9, "5-14",  19333,  9.4),
    @(2021, 9, "15-34", 102666, 6.3),
    @(2021, 9, "35-59", 156666, 6.1),
    @(2021, 9, "60-79", 88666,  5.3),
    @(2021, 9, ">=80",  45333,  5.3)
)
$r = 198
foreach ($row in $wk9) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 7. Move the frozen-pane scroll position / active selection the way a user
#    would end up after entering this data (cosmetic, best effort).
# ---------------------------------------------------------------------------
$ws.Range("E200").Select()
